# Solve Leetcode - 133. Clone Graph - DFS
# Adds a new entry (row 31) to the "Neetcode 150" tracker sheet for
# 695. Max Area of Island, mirroring the formatting of the preceding
# "200. Number of Islands" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$category = $ws.Range("A31")
$difficulty = $ws.Range("B31")
$name = $ws.Range("C31")
$notes = $ws.Range("D31")

$category.Value() = "Graphs"
$difficulty.Value() = "Medium"
$name.Value() = "695. Max Area of Island"
$notes.Value() = 'Whenever we encounter an island ("1"), run a dfs with all 4 directions to mark the entire island "#" for visited and track cur_max and max_max. Then go find the next island'

# Hyperlink the problem name cell to the LeetCode problem page (same
# pattern used by every other row in the sheet, where the hyperlink's
# display text is the target URL itself).
$ws.Hyperlinks.Add($name, "https://leetcode.com/problems/max-area-of-island/", "", "", "https://leetcode.com/problems/max-area-of-island/")

# Hyperlinks.Add() replaces the cell's text with TextToDisplay, so put
# the problem name back as the visible cell content.
$name.Value() = "695. Max Area of Island"

# Match the visual style used by row 30 (category/name use the built-in
# "Neutral"/"Good" cell styles; the notes column already inherits its
# wrapped/top-aligned style from the column default, same as D30).
$difficulty.Style() = "Neutral"
$name.Style() = "Good"

# Row 30 uses an auto-fit height of 28.8 for its two-line content; match
# it for the new row.
$ws.Rows.Item(31).RowHeight() = 28.8

# Move the active selection past the newly added row, as in the source
# workbook.
[void]$ws.Range("A32").Select()
